$d = $word.ActiveDocument

# --- 1) Split the "- Serve como base..." run so that a (now-current)
#        "_GoBack" edit-position bookmark sits right after "- Ser".
$find = $d.Content
$found = $find.Find.Execute("- Serve como base", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
if (-not $found) {
    throw "Could not locate target sentence"
}
$splitPoint = $find.Start + 5   # length of "- Ser"
$r = $d.Range($splitPoint, $splitPoint)
$r.Bookmarks.Add("_GoBack")

# --- 2) Remove the previous "_GoBack" bookmark that used to sit after
#        "Finally" in the Try/Catch/Finally table.
$bm = $d.Bookmarks
if ($bm.Exists("_GoBack")) {
    $old = $bm.Item("_GoBack")
    $old.Delete()
}
